$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right above the current row 1070, shifting all
# existing rows (old 1070-1129) down to 1072-1131.
$ws.Range("A1070:A1071").EntireRow.Insert()

# ---- New row 1070 ----
$ws.Cells.Item(1070, 1).Value  = 3
$ws.Cells.Item(1070, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(1070, 3).Value  = "Coquimbo"
$ws.Cells.Item(1070, 4).Value  = 44706
$ws.Cells.Item(1070, 5).Value  = 5
$ws.Cells.Item(1070, 6).Value  = 100112020
$ws.Cells.Item(1070, 7).Value  = "Tomate"
$ws.Cells.Item(1070, 8).Value  = "Larga vida"
$ws.Cells.Item(1070, 9).Value  = "Primera"
$ws.Cells.Item(1070, 10).Value = 480
$ws.Cells.Item(1070, 11).Value = 17000
$ws.Cells.Item(1070, 12).Value = 18000
$ws.Cells.Item(1070, 13).Value = 17521
$ws.Cells.Item(1070, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1070, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1070, 16).Value = 973
$ws.Cells.Item(1070, 17).Value = 18
$ws.Cells.Item(1070, 18).Value = "Hortaliza"

# ---- New row 1071 ----
$ws.Cells.Item(1071, 1).Value  = 3
$ws.Cells.Item(1071, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(1071, 3).Value  = "Coquimbo"
$ws.Cells.Item(1071, 4).Value  = 44706
$ws.Cells.Item(1071, 5).Value  = 5
$ws.Cells.Item(1071, 6).Value  = 100112020
$ws.Cells.Item(1071, 7).Value  = "Tomate"
$ws.Cells.Item(1071, 8).Value  = "Larga vida"
$ws.Cells.Item(1071, 9).Value  = "Segunda"
$ws.Cells.Item(1071, 10).Value = 180
$ws.Cells.Item(1071, 11).Value = 14000
$ws.Cells.Item(1071, 12).Value = 14000
$ws.Cells.Item(1071, 13).Value = 14000
$ws.Cells.Item(1071, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1071, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1071, 16).Value = 778
$ws.Cells.Item(1071, 17).Value = 18
$ws.Cells.Item(1071, 18).Value = "Hortaliza"
